$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '버스'
$ws.Range("C2").Value = 782
$ws.Range("B3").Value = '감사'
$ws.Range("C3").Value = 503
$ws.Range("B4").Value = '분'
$ws.Range("C4").Value = 466
$ws.Range("B5").Value = '기사님'
$ws.Range("C5").Value = 432
$ws.Range("B6").Value = '없'
$ws.Range("C6").Value = 362
$ws.Range("B7").Value = '글'
$ws.Range("C7").Value = 265
$ws.Range("B8").Value = '사람'
$ws.Range("C8").Value = 263
$ws.Range("B9").Value = '좋'
$ws.Range("C9").Value = 262
$ws.Range("B10").Value = '많'
$ws.Range("C10").Value = 222
$ws.Range("B11").Value = '통근'
$ws.Range("C11").Value = 220
$ws.Range("B12").Value = '보안'
$ws.Range("C12").Value = 220
$ws.Range("B13").Value = '이용'
$ws.Range("C13").Value = 217
$ws.Range("B14").Value = '생각'
$ws.Range("C14").Value = 192
$ws.Range("B15").Value = '문의'
$ws.Range("C15").Value = 177
$ws.Range("B16").Value = '차량'
$ws.Range("C16").Value = 177
$ws.Range("B17").Value = '셔틀'
$ws.Range("C17").Value = 176
$ws.Range("B18").Value = '불편'
$ws.Range("C18").Value = 173
$ws.Range("B19").Value = '관련'
$ws.Range("C19").Value = 170
$ws.Range("B20").Value = '확인'
$ws.Range("C20").Value = 164
$ws.Range("B21").Value = '퇴근'
$ws.Range("C21").Value = 160
$ws.Range("B22").Value = '칭찬'
$ws.Range("C22").Value = 157
$ws.Range("B23").Value = '직원'
$ws.Range("C23").Value = 156
$ws.Range("B24").Value = '부탁'
$ws.Range("C24").Value = 156
$ws.Range("B25").Value = '자리'
$ws.Range("C25").Value = 146
$ws.Range("B26").Value = '시간'
$ws.Range("C26").Value = 143
$ws.Range("B27").Value = '회사'
$ws.Range("C27").Value = 143
$ws.Range("B28").Value = '사원'
$ws.Range("C28").Value = 142
$ws.Range("B29").Value = '친절'
$ws.Range("C29").Value = 141
$ws.Range("B30").Value = '출근'
$ws.Range("C30").Value = 141
$ws.Range("B31").Value = '시'
$ws.Range("C31").Value = 141
$ws.Range("B32").Value = '사내'
$ws.Range("C32").Value = 135
$ws.Range("B33").Value = '하이닉스'
$ws.Range("C33").Value = 132
$ws.Range("B34").Value = '앞'
$ws.Range("C34").Value = 129
$ws.Range("B35").Value = '일'
$ws.Range("C35").Value = 126
$ws.Range("B36").Value = '말'
$ws.Range("C36").Value = 126
$ws.Range("B37").Value = '전'
$ws.Range("C37").Value = 125
$ws.Range("B38").Value = '번호'
$ws.Range("C38").Value = 124
$ws.Range("B39").Value = '안'
$ws.Range("C39").Value = 121
$ws.Range("B40").Value = '안녕'
$ws.Range("C40").Value = 119
$ws.Range("B41").Value = '노선'
$ws.Range("C41").Value = 118
$ws.Range("B42").Value = '사항'
$ws.Range("C42").Value = 116
$ws.Range("B43").Value = '관광'
$ws.Range("C43").Value = 113
$ws.Range("B44").Value = '좌석'
$ws.Range("C44").Value = 112
$ws.Range("B45").Value = '후'
$ws.Range("C45").Value = 108
$ws.Range("B46").Value = '답변'
$ws.Range("C46").Value = 107
$ws.Range("B47").Value = '하'
$ws.Range("C47").Value = 107
$ws.Range("B48").Value = '전화'
$ws.Range("C48").Value = 105
$ws.Range("B49").Value = '사용'
$ws.Range("C49").Value = 100
$ws.Range("B50").Value = '경우'
$ws.Range("C50").Value = 98
$ws.Range("B51").Value = '차'
$ws.Range("C51").Value = 98
$ws.Range("B52").Value = '분실물'
$ws.Range("C52").Value = 97
$ws.Range("B53").Value = '담당자'
$ws.Range("C53").Value = 97
$ws.Range("B54").Value = '옥계'
$ws.Range("C54").Value = 97
$ws.Range("B55").Value = '운행'
$ws.Range("C55").Value = 96
$ws.Range("B56").Value = '오늘'
$ws.Range("C56").Value = 96
$ws.Range("B57").Value = '운전'
$ws.Range("C57").Value = 95
$ws.Range("B58").Value = '중'
$ws.Range("C58").Value = 95
$ws.Range("B59").Value = '문'
$ws.Range("C59").Value = 93
$ws.Range("B60").Value = '주차장'
$ws.Range("C60").Value = 92
$ws.Range("B61").Value = '구리'
$ws.Range("C61").Value = 92
$ws.Range("B62").Value = '기분'
$ws.Range("C62").Value = 91
$ws.Range("B63").Value = '소리'
$ws.Range("C63").Value = 90
$ws.Range("B64").Value = '임산부'
$ws.Range("C64").Value = 89
$ws.Range("B65").Value = '태'
$ws.Range("C65").Value = 89
$ws.Range("B66").Value = '선'
$ws.Range("C66").Value = 88
$ws.Range("B67").Value = '어떻'
$ws.Range("C67").Value = 86
$ws.Range("B68").Value = '하세'
$ws.Range("C68").Value = 84
$ws.Range("B69").Value = '요청'
$ws.Range("C69").Value = 84
$ws.Range("B70").Value = '노'
$ws.Range("C70").Value = 84
$ws.Range("B71").Value = '관리'
$ws.Range("C71").Value = 83
$ws.Range("B72").Value = '안전'
$ws.Range("C72").Value = 83
$ws.Range("B73").Value = '인사'
$ws.Range("C73").Value = 82
$ws.Range("B74").Value = '수고'
$ws.Range("C74").Value = 82
$ws.Range("B75").Value = '연락'
$ws.Range("C75").Value = 82
$ws.Range("B76").Value = '로그'
$ws.Range("C76").Value = 82
$ws.Range("B77").Value = '서비스'
$ws.Range("C77").Value = 82
$ws.Range("B78").Value = '이렇'
$ws.Range("C78").Value = 81
$ws.Range("B79").Value = '말씀'
$ws.Range("C79").Value = 81
$ws.Range("B80").Value = '개선'
$ws.Range("C80").Value = 80
$ws.Range("B81").Value = '안녕하'
$ws.Range("C81").Value = 80
$ws.Range("B82").Value = '이스텍'
$ws.Range("C82").Value = 80
$ws.Range("B83").Value = '모습'
$ws.Range("C83").Value = 79
$ws.Range("B84").Value = '분실'
$ws.Range("C84").Value = 78
$ws.Range("B85").Value = '아침'
$ws.Range("C85").Value = 78
$ws.Range("B86").Value = '근무'
$ws.Range("C86").Value = 78
$ws.Range("B87").Value = '등'
$ws.Range("C87").Value = 78
$ws.Range("B88").Value = '식당'
$ws.Range("C88").Value = 77
$ws.Range("B89").Value = '요원'
$ws.Range("C89").Value = 77
$ws.Range("B90").Value = '제가'
$ws.Range("C90").Value = 76
$ws.Range("B91").Value = '문제'
$ws.Range("C91").Value = 75
$ws.Range("B92").Value = '인'
$ws.Range("C92").Value = 74
$ws.Range("B93").Value = '가능'
$ws.Range("C93").Value = 73
$ws.Range("B94").Value = '구성원'
$ws.Range("C94").Value = 73
$ws.Range("B95").Value = '그렇'
$ws.Range("C95").Value = 73
$ws.Range("B96").Value = '주시'
$ws.Range("C96").Value = 72
$ws.Range("B97").Value = '뒤'
$ws.Range("C97").Value = 72
$ws.Range("B98").Value = '이천'
$ws.Range("C98").Value = 71
$ws.Range("B99").Value = '출퇴근'
$ws.Range("C99").Value = 70
$ws.Range("B100").Value = '마음'
$ws.Range("C100").Value = 69
$ws.Range("B101").Value = '때문'
$ws.Range("C101").Value = 69
